$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "pfm" column (G) to sit right after "pylddt_model" (B), i.e. before the
# two columns that are about to be removed.
$ws.Range("G:G").Cut() | Out-Null
$ws.Range("C:C").Insert() | Out-Null

# Remove the now-obsolete "pylddt(or)" / "pylddt(and)" columns (shifted to D:E).
$ws.Range("D:E").Delete() | Out-Null

# Rename headers.
$ws.Range("C1").Value = "pylddt"
$ws.Range("A1").Value = "set"

# Match the number format used by the rest of the decimal columns.
$ws.Range("C2:C7").NumberFormat = "0.000"
$ws.Range("C2:C7").HorizontalAlignment = -4108

# Bold header row.
$ws.Range("A1:E1").Font.Bold = $true

# Restore the expected column widths (the cut/insert above drags each
# column's original width along with it).
$ws.Range("B:C").ColumnWidth = 14.28515625
$ws.Range("D:D").ColumnWidth = 15
$ws.Range("E:E").ColumnWidth = 9.140625

